$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix bug description in B23 and add a new note in C23
$ws.Range("B23").Value = "не удаляет и не правит событие"
$ws.Range("C23").Value = "почему-то event после того,как форма открывается, обнуляется undefined"

# Move the active selection to A23
$ws.Range("A23").Select()

$wb.Save()
